# Applies: "Added distance formula to point."
#
# Summary of the change:
#  - A new worksheet "Sheet3" is appended after Sheet2, becomes the active
#    (selected) sheet/tab.
#  - Sheet3 contains two 3-D points (row 1 and row 2) and, in row 4, the
#    squared per-axis differences plus the overall Euclidean distance
#    (square root of the sum of squares) in E4.
#  - The view/selection on Sheet1 moves from F32 to B42 (scrolled so A13 is
#    the top-left visible cell).
#  - Sheet2 is no longer the active tab (Sheet3 is now active instead).

$wb = $excel.ActiveWorkbook

# --- Sheet1: update the view/selection (Sheet1 is no longer the active tab
#     either, but its scroll position & selection change as in the diff) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 13
$win.ScrollColumn = 1
$ws1.Range("B42").Select()

# --- Sheet2: no cell-data changes; it simply stops being the active tab
#     once Sheet3 is added & activated below. ---

# --- Sheet3: brand-new worksheet appended at the end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)
$ws3.Name = "Sheet3"

# Point 1
$ws3.Range("A1").Value = 1
$ws3.Range("B1").Value = 2
$ws3.Range("C1").Value = 3

# Point 2
$ws3.Range("A2").Value = 5
$ws3.Range("B2").Value = 8.5
$ws3.Range("C2").Value = 9.3

# Squared differences per axis
$ws3.Range("A4").Formula = "=(A2-A1)^2"
$ws3.Range("B4").Formula = "=(B2-B1)^2"
$ws3.Range("C4").Formula = "=(C2-C1)^2"

# Euclidean distance between the two points
$ws3.Range("E4").Formula = "=SQRT(SUM(A4:C4))"

# Leave the selection/cursor on the cell right after the formula, as in the
# authored workbook.
$ws3.Range("E5").Select()
